$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated symbol list refresh: new prices/volumes/hour stamp for each coin row,
# plus three rows (17-24) that shifted position with new Coin/Link pairs.
foreach ($addr in @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "D19", "E19", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "D24", "E24", "G24", "E25", "G25", "D26", "E26", "G26", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "D38", "E38", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46", "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50", "G50", "D51", "E51", "G51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "321.67"
$ws.Range("E2").Value = "-3.07%"
$ws.Range("G2").Value = "7"
$ws.Range("D3").Value = "43.08"
$ws.Range("E3").Value = "-5.88%"
$ws.Range("G3").Value = "7"
$ws.Range("D4").Value = "5.198"
$ws.Range("E4").Value = "-7.22%"
$ws.Range("G4").Value = "7"
$ws.Range("D5").Value = "0.08193"
$ws.Range("E5").Value = "-1.92%"
$ws.Range("G5").Value = "7"
$ws.Range("D6").Value = "4.324"
$ws.Range("E6").Value = "-3.43%"
$ws.Range("G6").Value = "7"
$ws.Range("D7").Value = "1.829"
$ws.Range("E7").Value = "-10.97%"
$ws.Range("G7").Value = "7"
$ws.Range("D8").Value = "0.9387"
$ws.Range("E8").Value = "-3.86%"
$ws.Range("G8").Value = "7"
$ws.Range("D9").Value = "0.1115"
$ws.Range("E9").Value = "-3.99%"
$ws.Range("G9").Value = "7"
$ws.Range("D10").Value = "0.1859"
$ws.Range("E10").Value = "-3.28%"
$ws.Range("G10").Value = "7"
$ws.Range("D11").Value = "0.09379"
$ws.Range("E11").Value = "-5.45%"
$ws.Range("G11").Value = "7"
$ws.Range("D12").Value = "0.04619"
$ws.Range("E12").Value = "-1.18%"
$ws.Range("G12").Value = "7"
$ws.Range("D13").Value = "7.416"
$ws.Range("E13").Value = "-28.60%"
$ws.Range("G13").Value = "7"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").Value = "0.01%"
$ws.Range("G14").Value = "7"
$ws.Range("D15").Value = "0.001279"
$ws.Range("E15").Value = "-0.30%"
$ws.Range("G15").Value = "7"
$ws.Range("D16").Value = "0.005937"
$ws.Range("E16").Value = "-1.51%"
$ws.Range("G16").Value = "7"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.357"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("G17").Value = "7"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.514"
$ws.Range("E18").Value = "-2.25%"
$ws.Range("G18").Value = "7"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3366"
$ws.Range("E19").Value = "0.09%"
$ws.Range("G19").Value = "7"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1388"
$ws.Range("E20").Value = "-0.97%"
$ws.Range("G20").Value = "7"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.2620"
$ws.Range("E21").Value = "-1.21%"
$ws.Range("G21").Value = "7"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04158"
$ws.Range("E22").Value = "-0.84%"
$ws.Range("G22").Value = "7"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001250"
$ws.Range("E23").Value = "-4.69%"
$ws.Range("G23").Value = "7"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.004307"
$ws.Range("E24").Value = "-7.37%"
$ws.Range("G24").Value = "7"
$ws.Range("E25").Value = "-14.32%"
$ws.Range("G25").Value = "7"
$ws.Range("D26").Value = "0.0002978"
$ws.Range("E26").Value = "-20.51%"
$ws.Range("G26").Value = "7"
$ws.Range("G27").Value = "7"
$ws.Range("G28").Value = "7"
$ws.Range("G29").Value = "7"
$ws.Range("G30").Value = "7"
$ws.Range("G31").Value = "7"
$ws.Range("G32").Value = "7"
$ws.Range("G33").Value = "7"
$ws.Range("G34").Value = "7"
$ws.Range("G35").Value = "7"
$ws.Range("G36").Value = "7"
$ws.Range("G37").Value = "7"
$ws.Range("D38").Value = "0.02723"
$ws.Range("E38").Value = "-1.36%"
$ws.Range("G38").Value = "7"
$ws.Range("D39").Value = "0.05547"
$ws.Range("E39").Value = "-4.34%"
$ws.Range("G39").Value = "7"
$ws.Range("D40").Value = "0.007973"
$ws.Range("E40").Value = "3.10%"
$ws.Range("G40").Value = "7"
$ws.Range("D41").Value = "0.1395"
$ws.Range("E41").Value = "-2.85%"
$ws.Range("G41").Value = "7"
$ws.Range("D42").Value = "0.006540"
$ws.Range("E42").Value = "-10.03%"
$ws.Range("G42").Value = "7"
$ws.Range("D43").Value = "0.002098"
$ws.Range("E43").Value = "4.15%"
$ws.Range("G43").Value = "7"
$ws.Range("D44").Value = "0.008228"
$ws.Range("E44").Value = "1.48%"
$ws.Range("G44").Value = "7"
$ws.Range("D45").Value = "0.3201"
$ws.Range("E45").Value = "-5.99%"
$ws.Range("G45").Value = "7"
$ws.Range("D46").Value = "0.00006954"
$ws.Range("E46").Value = "-4.78%"
$ws.Range("G46").Value = "7"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("G47").Value = "7"
$ws.Range("D48").Value = "0.003464"
$ws.Range("E48").Value = "-1.02%"
$ws.Range("G48").Value = "7"
$ws.Range("D49").Value = "0.003529"
$ws.Range("E49").Value = "0.71%"
$ws.Range("G49").Value = "7"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("G50").Value = "7"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.21%"
$ws.Range("G51").Value = "7"

# Drop the temporary Text number format so the cells keep the workbook's
# default (unstyled) appearance, matching the source formatting.
foreach ($addr in @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "D19", "E19", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "D24", "E24", "G24", "E25", "G25", "D26", "E26", "G26", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "D38", "E38", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46", "G46", "D47", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50", "G50", "D51", "E51", "G51")) {
    $ws.Range($addr).Style = "Normal"
}
